# "redondeado a 4 cifras" - update the DemandaUPME projection table with the
# more precise (4-decimal) figures and switch the number format from
# "0.000" (3 decimals) to "0.0000" (4 decimals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DemandaUPME")

# --- Updated data values (columns B-G, rows 3-17) ------------------------
$ws.Range("B3").Value = 11184.005999999999
$ws.Range("C3").Value = 10800.5607
$ws.Range("D3").Value = 10430.642900000001
$ws.Range("E3").Value = 0.072099999999999997
$ws.Range("F3").Value = 0.035299999999999998
$ws.Range("G3").Value = -0.0001

$ws.Range("B4").Value = 11485.919599999999
$ws.Range("C4").Value = 11095.8215
$ws.Range("D4").Value = 10719.485699999999
$ws.Range("E4").Value = 0.027
$ws.Range("F4").Value = 0.027300000000000001
$ws.Range("G4").Value = 0.027699999999999999

$ws.Range("B5").Value = 11738.5939
$ws.Range("C5").Value = 11343.4645
$ws.Range("D5").Value = 10962.2749
$ws.Range("E5").Value = 0.021999999999999999
$ws.Range("F5").Value = 0.0223
$ws.Range("G5").Value = 0.022599999999999999

$ws.Range("B6").Value = 11971.9085
$ws.Range("C6").Value = 11571.3709
$ws.Range("D6").Value = 11184.963900000001
$ws.Range("E6").Value = 0.019900000000000001
$ws.Range("F6").Value = 0.0201
$ws.Range("G6").Value = 0.020299999999999999

$ws.Range("B7").Value = 12196.2287
$ws.Range("C7").Value = 11790.809800000001
$ws.Range("D7").Value = 11399.6937
$ws.Range("E7").Value = 0.018700000000000001
$ws.Range("F7").Value = 0.019
$ws.Range("G7").Value = 0.019199999999999998

$ws.Range("B8").Value = 12393.3611
$ws.Range("C8").Value = 11980.0988
$ws.Range("D8").Value = 11581.4161
$ws.Range("E8").Value = 0.016199999999999999
$ws.Range("F8").Value = 0.0161
$ws.Range("G8").Value = 0.015900000000000001

$ws.Range("B9").Value = 12528.194
$ws.Range("C9").Value = 12109.5229
$ws.Range("D9").Value = 11705.622100000001
$ws.Range("E9").Value = 0.0109
$ws.Range("F9").Value = 0.010800000000000001
$ws.Range("G9").Value = 0.010699999999999999

$ws.Range("B10").Value = 12756.818300000001
$ws.Range("C10").Value = 12333.4049
$ws.Range("D10").Value = 11924.9292
$ws.Range("E10").Value = 0.018200000000000001
$ws.Range("F10").Value = 0.018499999999999999
$ws.Range("G10").Value = 0.018700000000000001

$ws.Range("B11").Value = 12892.9167
$ws.Range("C11").Value = 12464.799300000001
$ws.Range("D11").Value = 12051.785400000001
$ws.Range("E11").Value = 0.010699999999999999
$ws.Range("F11").Value = 0.010699999999999999
$ws.Range("G11").Value = 0.0106

$ws.Range("B12").Value = 13043.571900000001
$ws.Range("C12").Value = 12610.6196
$ws.Range("D12").Value = 12192.9414
$ws.Range("E12").Value = 0.0117
$ws.Range("F12").Value = 0.0117
$ws.Range("G12").Value = 0.0117

$ws.Range("B13").Value = 13280.8086
$ws.Range("C13").Value = 12839.5717
$ws.Range("D13").Value = 12413.901099999999
$ws.Range("E13").Value = 0.018200000000000001
$ws.Range("F13").Value = 0.018200000000000001
$ws.Range("G13").Value = 0.018100000000000002

$ws.Range("B14").Value = 13449.2245
$ws.Range("C14").Value = 13002.5852
$ws.Range("D14").Value = 12571.702799999999
$ws.Range("E14").Value = 0.012699999999999999
$ws.Range("F14").Value = 0.012699999999999999
$ws.Range("G14").Value = 0.012699999999999999

$ws.Range("B15").Value = 13608.7202
$ws.Range("C15").Value = 13156.831899999999
$ws.Range("D15").Value = 12720.8858
$ws.Range("E15").Value = 0.011900000000000001
$ws.Range("F15").Value = 0.011900000000000001
$ws.Range("G15").Value = 0.011900000000000001

$ws.Range("B16").Value = 13767.669400000001
$ws.Range("C16").Value = 13310.4512
$ws.Range("D16").Value = 12869.363300000001
$ws.Range("E16").Value = 0.0117
$ws.Range("F16").Value = 0.0117
$ws.Range("G16").Value = 0.0117

$ws.Range("B17").Value = 13931.0592
$ws.Range("C17").Value = 13468.280199999999
$ws.Range("D17").Value = 13021.827499999999
$ws.Range("E17").Value = 0.011900000000000001
$ws.Range("F17").Value = 0.011900000000000001
$ws.Range("G17").Value = 0.0118

# --- Number format: 3 decimals -> 4 decimals ------------------------------
$ws.Range("B3:G17").NumberFormat = "0.0000"

# --- Selection moves from J1:O1048576 to I3 -------------------------------
$ws.Activate() | Out-Null
$ws.Range("I3").Select() | Out-Null
